$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Experimental value (was blank) -> "false" (literal text, not boolean)
$ws.Cells.Item(7, 2).Formula = "=""false"""
$ws.Cells.Item(7, 2).Copy()
$ws.Cells.Item(7, 2).PasteSpecial(-4163)

# Date value updated
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"

# Description value (was blank) -> new description text
$ws.Range("B17").Value = "Cardiovascular risk categories based on fitness assessment"
